# This script reproduces the commit:
#   "fixed DeepL and started function-ifying"
#
# The underlying semantic change is: on the 'ZH test' sheet, several Chinese
# translation strings had trailing whitespace (regular spaces and/or
# non-breaking spaces) left over from a DeepL export. This script strips
# that trailing whitespace. The 'source' sheet pulls these values in via
# formulas, so its cached values update automatically on recalculation.

$wb = $excel.ActiveWorkbook
$zh = $wb.Worksheets.Item("ZH test")

# Trimmed (DeepL artifact whitespace removed) versions of A1:A19.
$values = @(
    "你的呼吸道感染是由冠状病毒引起的。",
    "犬冠状病毒病。",
    "新生犊腹泻冠状病毒。",
    "猪传染性胃肠炎冠状病毒。",
    "大熊猫犬冠状病毒的分离与鉴定。",
    "冠状病毒感染引起水貂肠炎的诊断。",
    "禽冠状病毒的介电松弛效应及灭活。",
    "冠状病毒基因组发明专利保护的思考。",
    "中和试验进一步证明该分离病毒是一种犬冠状病毒。",
    "基于密码子使用模式的冠状病毒亲缘关系分析。",
    "猪呼吸道冠状病毒及实验室诊断方法研究进展。",
    "与犬细小病毒，但没有具体的治疗犬冠状病毒。",
    "突起蛋白是冠状病毒的主要抗原，包含许多抗原决定簇。",
    "从健康狐狸、貉粪中检出犬冠状病毒的两种基因型。",
    "大熊猫犬冠状病毒部分纤突蛋白基因的扩增与序列分析。",
    "病毒被反病毒软件删除了。",
    "他染上了一种致命病毒。",
    "血液检验显示有病毒存在。",
    "哎呀，电脑又有病毒了。"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $zh.Range("A$row").Value = $values[$i]
    $zh.Range("A" + ($row + 19)).Value = $values[$i]
}

$wb.RefreshAll()
$excel.CalculateFullRebuild()

# Restore sheet view scroll/selection state as in the edit.
$src = $wb.Worksheets.Item("source")
$srcView = $src.Application.ActiveWindow
$src.Activate()
$src.Range("A1").Select()

$zh.Activate()
$zhWindow = $excel.ActiveWindow
$zhWindow.ScrollRow = 10
$zh.Range("I23").Select()

$src.Activate()
